$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append new row 29 with the latest test-mail entry ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A29").Value = "Wanneer zijn jullie open?"
$ws.Range("B29").Value = "mailmind.test@zohomail.eu"
$ws.Range("C29").Value = "Testmail #1: Wanneer zijn jullie open?"
$ws.Range("D29").Value = "Openingstijden / Locatie"
$ws.Range("F29").Value = "2025-06-26 22:57:10"
$ws.Range("G29").Value = "Nee"
$ws.Range("H29").Value = "Nee"
$ws.Range("I29").Value = "Nee"

# Extend the conditional-formatting ranges so they keep covering the new row
$fcD = $ws.Range("D2:D28").FormatConditions
$fcD.Item(1).ModifyAppliesToRange($ws.Range("D2:D29"))

$fcG = $ws.Range("G2:G28").FormatConditions
$fcG.Item(1).ModifyAppliesToRange($ws.Range("G2:G29"))

$fcH = $ws.Range("H2:H28").FormatConditions
$fcH.Item(1).ModifyAppliesToRange($ws.Range("H2:H29"))

$fcI = $ws.Range("I2:I28").FormatConditions
$fcI.Item(1).ModifyAppliesToRange($ws.Range("I2:I29"))

# --- "Dashboard" sheet: bump the "Openingstijden / Locatie" count ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 6
